$wb = $excel.ActiveWorkbook

# --- Sheet "Cycle_2021-2022": deselect tab, move selection from H31 to H37 ---
$ws2 = $wb.Worksheets.Item("Cycle_2021-2022")
$ws2.Activate()
$ws2.Range("H37").Select()

# --- Sheet "Cycle_2020-2021 argv": fill in prediction / delta formulas on rows 8-11 ---
$ws3 = $wb.Worksheets.Item("Cycle_2020-2021 argv")

# Row 8: predicted avg GPA of added
$ws3.Range("E8").Formula = "=(E6-D6*(D2/E2))/(E10/E2)"
$ws3.Range("F8").Formula = "=(F6-E6*(E2/F2))/(F10/F2)"
$ws3.Range("G8").Formula = "=(G6-F6*(F2/G2))/(G10/G2)"
$ws3.Range("H8").Formula = "=(H6-G6*(G2/H2))/(H10/H2)"
$ws3.Range("I8").Formula = "=(I6-H6*(H2/I2))/(I10/I2)"
$ws3.Range("J8").Formula = "=(J6-I6*(I2/J2))/(J10/J2)"
$ws3.Range("K8").Formula = "=(K6-J6*(J2/K2))/(K10/K2)"
$ws3.Range("M8").Formula = "=(M6-K6*(K2/M2))/(M10/M2)"
$ws3.Range("N8").Formula = "=(N6-M6*(M2/N2))/(N10/N2)"
$ws3.Range("O8").Formula = "=(O6-N6*(N2/O2))/(O10/O2)"
$ws3.Range("P8").Formula = "=(P6-O6*(O2/P2))/(P10/P2)"
$ws3.Range("Q8").Formula = "=(Q6-P6*(P2/Q2))/(Q10/Q2)"
$ws3.Range("R8").Formula = "=(R6-Q6*(Q2/R2))/(R10/R2)"
$ws3.Range("S8").Formula = "=(S6-R6*(R2/S2))/(S10/S2)"

# Row 9: predicted average MCAT of added
$ws3.Range("E9").Formula = "=(E7-D7*(D2/E2))/(E10/E2)"
$ws3.Range("F9").Formula = "=(F7-E7*(E2/F2))/(F10/F2)"
$ws3.Range("G9").Formula = "=(G7-F7*(F2/G2))/(G10/G2)"
$ws3.Range("H9").Formula = "=(H7-G7*(G2/H2))/(H10/H2)"
$ws3.Range("I9").Formula = "=(I7-H7*(H2/I2))/(I10/I2)"
$ws3.Range("J9").Formula = "=(J7-I7*(I2/J2))/(J10/J2)"
$ws3.Range("K9").Formula = "=(K7-J7*(J2/K2))/(K10/K2)"
$ws3.Range("M9").Formula = "=(M7-K7*(K2/M2))/(M10/M2)"
$ws3.Range("N9").Formula = "=(N7-M7*(M2/N2))/(N10/N2)"
$ws3.Range("O9").Formula = "=(O7-N7*(N2/O2))/(O10/O2)"
$ws3.Range("P9").Formula = "=(P7-O7*(O2/P2))/(P10/P2)"
$ws3.Range("Q9").Formula = "=(Q7-P7*(P2/Q2))/(Q10/Q2)"
$ws3.Range("R9").Formula = "=(R7-Q7*(Q2/R2))/(R10/R2)"
$ws3.Range("S9").Formula = "=(S7-R7*(R2/S2))/(S10/S2)"

# Row 10: # Added
$ws3.Range("E10").Formula = "=E2-D2"
$ws3.Range("F10").Formula = "=F2-E2"
$ws3.Range("G10").Formula = "=G2-F2"
$ws3.Range("H10").Formula = "=H2-G2"
$ws3.Range("I10").Formula = "=I2-H2"
$ws3.Range("J10").Formula = "=J2-I2"
$ws3.Range("K10").Formula = "=K2-J2"
$ws3.Range("L10").Formula = "=L2-K2"
$ws3.Range("M10").Formula = "=M2-K2"
$ws3.Range("N10").Formula = "=N2-M2"
$ws3.Range("O10").Formula = "=O2-N2"
$ws3.Range("P10").Formula = "=P2-O2"
$ws3.Range("Q10").Formula = "=Q2-P2"
$ws3.Range("R10").Formula = "=R2-Q2"
$ws3.Range("S10").Formula = "=S2-R2"

# Row 11: # Rs reversed
$ws3.Range("E11").Formula = "=ABS(E4-C4)"
$ws3.Range("F11").Formula = "=ABS(F4-E4)"
$ws3.Range("G11").Formula = "=ABS(G4-F4)"
$ws3.Range("H11").Formula = "=ABS(H4-G4)"
$ws3.Range("I11").Formula = "=ABS(I4-H4)"
$ws3.Range("J11").Formula = "=ABS(J4-I4)"
$ws3.Range("K11").Formula = "=ABS(K4-J4)"
$ws3.Range("M11").Formula = "=ABS(M4-K4)"
$ws3.Range("N11").Formula = "=ABS(N4-M4)"
$ws3.Range("O11").Formula = "=ABS(O4-N4)"
$ws3.Range("P11").Formula = "=ABS(P4-O4)"
$ws3.Range("Q11").Formula = "=ABS(Q4-P4)"

$wb.Application.CalculateFull()

# Activate this sheet last so it becomes the selected tab / active sheet in the workbook
$ws3.Activate()
$ws3.Range("L20").Select()
